# quickstart-images.pptx — "update timeservice interface in quickstart tutorial"
#
# Slide 5 ("Textfeld 40"): the "+ getLocation()" / "+ subscribe()" method rows
# get split into three runs apiece (prefix / identifier / parens) so the
# identifier can be flagged independently (spell-check styling in the
# original authoring session). Slide 6 ("Text Box 12"): the "subscribe()"
# call-out is split into "subscribe" + "(...)" using an ellipsis for the
# argument list.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5 — Textfeld 40: "+ getLocation()" / "+ subscribe()"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(4)
$tr5 = $sh5.TextFrame.TextRange

# Paragraph 1: "+ getLocation()" -> "+ " / "getLocation" / "()"
$tr5.Characters(1, 2).Text = "+ "
$tr5.Characters(3, 11).Text = "getLocation"
$tr5.Characters(14, 2).Text = "()"

# Paragraph 2: "+ subscribe()" -> "+ " / "subscribe" / "()"
$tr5.Characters(17, 2).Text = "+ "
$tr5.Characters(19, 9).Text = "subscribe"
$tr5.Characters(28, 2).Text = "()"

# ---------------------------------------------------------------------
# Slide 6 — Text Box 12: "subscribe()" -> "subscribe" + "(...)"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(4)
$tr6 = $sh6.TextFrame.TextRange

# Split off the trailing "()" from "subscribe" first (re-assigning the same
# text forces the run boundary without touching formatting), then swap the
# parens' content for an ellipsis placeholder.
$tr6.Characters(1, 9).Text = "subscribe"
$tr6.Characters(10, 2).Text = "(" + [char]0x2026 + ")"
